{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is the async (context) => { ... } function.\n\n// List of exact text replacements to perform. Order matters only in that\n// each pair is matched by an exact, whole-string search, so distinct source\n// strings never collide with each other (even when old/new values are\n// \"swapped\" between two bullets).\nconst replacements = [\n  {\n    oldText: \"Play Dragon Egg Free Slot Review | Exciting Free Spins Mode\",\n    newText: \"Play Dragon Egg Free Online Slot\"\n  },\n  {\n    oldText: \"Potentially substantial winnings in free spins mode\",\n    newText: \"Slightly higher value wins than average\"\n  },\n  {\n    oldText: \"Clear and concise graphics\",\n    newText: \"Simple and clear graphics and symbols\"\n  },\n  {\n    oldText: \"Well-suited sound and music\",\n    newText: \"Background music fits the game's setting\"\n  },\n  {\n    oldText: \"Higher value wins than average\",\n    newText: \"Potential for substantial winnings in free spins mode\"\n  },\n  {\n    oldText: \"Difficult to trigger free spins mode\",\n    newText: \"Difficult to trigger the free spins mode\"\n  },\n  {\n    oldText: \"Limited symbol variety\",\n    newText: \"Limited variety in dragon symbol designs\"\n  },\n  {\n    oldText:\n      \"Read our review of Dragon Egg, an online slot game with a free spins mode. Enjoy higher value wins than average and clear graphics. Play Dragon Egg for free.\",\n    newText:\n      \"Read our review of Dragon Egg and play this free online slot game with substantial winnings.\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# wdFindContinue = 1, wdReplaceAll = 2 (standard Word enum values).\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n\n# Title heading (also repeated later in bold near the end of the document).\nReplace-Text \"Play Dragon Egg Free Slot Review | Exciting Free Spins Mode\" \"Play Dragon Egg Free Online Slot\"\n\n# \"What we like\" bullet list.\nReplace-Text \"Potentially substantial winnings in free spins mode\" \"Slightly higher value wins than average\"\nReplace-Text \"Clear and concise graphics\" \"Simple and clear graphics and symbols\"\nReplace-Text \"Well-suited sound and music\" \"Background music fits the game's setting\"\nReplace-Text \"Higher value wins than average\" \"Potential for substantial winnings in free spins mode\"\n\n# \"What we don't like\" bullet list.\nReplace-Text \"Difficult to trigger free spins mode\" \"Difficult to trigger the free spins mode\"\nReplace-Text \"Limited symbol variety\" \"Limited variety in dragon symbol designs\"\n\n# Closing meta title/description paragraphs.\nReplace-Text \"Read our review of Dragon Egg, an online slot game with a free spins mode. Enjoy higher value wins than average and clear graphics. Play Dragon Egg for free.\" \"Read our review of Dragon Egg and play this free online slot game with substantial winnings.\"\n"}
